# Add a new "IS_CLICKABLE" test step between the existing NAVIGATE_TO step
# and the CLICK step on the TestCases sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

# Insert a new row above row 3 - this pushes rows 3-8 down to 4-9 and
# keeps their formatting/styles intact (Excel copies formatting from the row above).
$ws.Rows.Item(3).EntireRow.Insert()

# Fill in the new row 3 with the IS_CLICKABLE step.
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = "Check if Today's Deals link is clickable"
$ws.Range("F3").Value = "IS_CLICKABLE"
$ws.Range("H3").Value = "xpath"
$ws.Range("I3").Value = ".//*[@id='nav-xshop']/a[2]"
$ws.Range("I3").Style = $ws.Range("I4").Style

# Renumber the testStepID column for the steps that got pushed down.
$ws.Range("D4").Value = 3
$ws.Range("D5").Value = 4
$ws.Range("D6").Value = 5

# Clear out the stray leftover xpath values that used to live in row 7
# and now sit in row 8 after the insert (they are being removed).
$ws.Range("H8").ClearContents()
$ws.Range("I8").ClearContents()

# Restore the selection to A3 (matches where the author was working).
$ws.Range("A3").Select()
